{"js": "const body = context.document.body;\n\n// Each (old, new) pair below corresponds to one <w:t> run that changed in the diff.\n// All old values are unique in the document, so a literal, case-sensitive,\n// whole-match search safely identifies the single run to update.\nconst replacements = [\n  [\"2024-09-25 Wednesday\", \"2024-09-26 Thursday\"],\n  [\"88-40=\", \"87-73=\"],\n  [\"39-9=\", \"39-33=\"],\n  [\"55+44=\", \"32+25=\"],\n  [\"81-13=\", \"76-7=\"],\n  [\"71+3=\", \"32-24=\"],\n  [\"61+7=\", \"78-33=\"],\n  [\"3+82=\", \"45+43=\"],\n  [\"27+23=\", \"50+9=\"],\n  [\"73+8=\", \"33+65=\"],\n  [\"78+20=\", \"70-59=\"],\n  [\"57-15=\", \"49+15=\"],\n  [\"13-7=\", \"11+17=\"],\n  [\"30+23=\", \"36-16=\"],\n  [\"31-0=\", \"0+40=\"],\n  [\"23+64=\", \"14+13=\"],\n  [\"58+10=\", \"34-18=\"],\n  [\"1+20=\", \"29+26=\"],\n  [\"29+52=\", \"5+66=\"],\n  [\"53+37=\", \"71-68=\"],\n  [\"35+18=\", \"19+7=\"],\n  [\"32+59=\", \"48-44=\"],\n  [\"35+63=\", \"36-5=\"],\n  [\"41+56=\", \"38+9=\"],\n  [\"79-22=\", \"56+15=\"],\n  [\"64+33=\", \"23+12=\"],\n  [\"22+65=\", \"69-5=\"],\n  [\"18-1=\", \"35-8=\"],\n  [\"6+29=\", \"39-17=\"],\n  [\"14-11=\", \"75-43=\"],\n  [\"76-1=\", \"86-47=\"],\n  [\"14+52=\", \"88-51=\"],\n  [\"43-12=\", \"38+5=\"],\n  [\"51+36=\", \"32+17=\"],\n  [\"4+85=\", \"93-79=\"],\n  [\"95-59=\", \"83-78=\"],\n  [\"73-21=\", \"95-37=\"],\n  [\"65+11=\", \"67+5=\"],\n  [\"14+15=\", \"38-21=\"],\n  [\"96-2=\", \"22+46=\"],\n  [\"92-45=\", \"61-41=\"],\n  [\"44+29=\", \"82-15=\"],\n  [\"56+35=\", \"14+44=\"],\n  [\"76+2=\", \"69-64=\"],\n  [\"39-24=\", \"44+5=\"],\n  [\"71-31=\", \"46+53=\"],\n  [\"21+15=\", \"45-21=\"],\n  [\"15+26=\", \"55+29=\"],\n  [\"64-43=\", \"63+15=\"],\n  [\"31+8=\", \"25+9=\"],\n  [\"99-87=\", \"86-35=\"],\n  [\"93+6=\", \"66-14=\"],\n  [\"67+31=\", \"79+9=\"],\n  [\"12+36=\", \"89-79=\"],\n  [\"83-62=\", \"31+33=\"],\n  [\"80-17=\", \"32+19=\"],\n  [\"70-44=\", \"79-56=\"],\n  [\"59-22=\", \"79-1=\"],\n  [\"49-13=\", \"77+14=\"],\n  [\"12+65=\", \"83-9=\"],\n  [\"47+43=\", \"42+47=\"],\n  [\"79-33=\", \"37+5=\"],\n  [\"33-2=\", \"6+66=\"],\n  [\"56-36=\", \"17-1=\"],\n  [\"67-53=\", \"3+84=\"],\n  [\"36-23=\", \"35+16=\"],\n  [\"64+31=\", \"85-56=\"],\n  [\"1+57=\", \"74-59=\"],\n  [\"58-22=\", \"10+88=\"],\n  [\"39-4=\", \"36-31=\"],\n  [\"29-18=\", \"88-68=\"],\n  [\"50-32=\", \"34-27=\"],\n  [\"38-38=\", \"80-34=\"],\n  [\"12+29=\", \"45-19=\"],\n  [\"83-44=\", \"51+43=\"],\n  [\"65-29=\", \"47-47=\"],\n  [\"96-56=\", \"98-86=\"],\n  [\"65-36=\", \"81-50=\"],\n  [\"33-21=\", \"83-21=\"],\n  [\"25+49=\", \"87-25=\"],\n  [\"83-52=\", \"82-68=\"],\n  [\"43+0=\", \"26+1=\"],\n  [\"38-4=\", \"40-27=\"],\n  [\"99-31=\", \"92-28=\"],\n  [\"24-22=\", \"38-35=\"],\n  [\"79-15=\", \"30-5=\"],\n  [\"42+44=\", \"2+94=\"],\n  [\"69+22=\", \"32+12=\"],\n  [\"12+77=\", \"66-19=\"],\n  [\"75-54=\", \"85-18=\"],\n  [\"52-35=\", \"13+51=\"],\n  [\"99-98=\", \"65+32=\"],\n  [\"84-2=\", \"41+52=\"],\n  [\"26+11=\", \"63+12=\"],\n  [\"73+0=\", \"56-44=\"],\n  [\"53-1=\", \"46+20=\"],\n  [\"76-72=\", \"47-19=\"],\n  [\"95-42=\", \"57+24=\"],\n  [\"6+91=\", \"81-40=\"],\n  [\"60+27=\", \"21+25=\"],\n  [\"41-9=\", \"9+72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n# Each (old, new) pair below corresponds to one run of text that changed in the diff.\n# All old values are unique in the document, so a literal, case-sensitive,\n# find-and-replace-all safely targets exactly the one run that should change.\n$replacements = @(\n    @('2024-09-25 Wednesday', '2024-09-26 Thursday'),\n    @('88-40=', '87-73='),\n    @('39-9=', '39-33='),\n    @('55+44=', '32+25='),\n    @('81-13=', '76-7='),\n    @('71+3=', '32-24='),\n    @('61+7=', '78-33='),\n    @('3+82=', '45+43='),\n    @('27+23=', '50+9='),\n    @('73+8=', '33+65='),\n    @('78+20=', '70-59='),\n    @('57-15=', '49+15='),\n    @('13-7=', '11+17='),\n    @('30+23=', '36-16='),\n    @('31-0=', '0+40='),\n    @('23+64=', '14+13='),\n    @('58+10=', '34-18='),\n    @('1+20=', '29+26='),\n    @('29+52=', '5+66='),\n    @('53+37=', '71-68='),\n    @('35+18=', '19+7='),\n    @('32+59=', '48-44='),\n    @('35+63=', '36-5='),\n    @('41+56=', '38+9='),\n    @('79-22=', '56+15='),\n    @('64+33=', '23+12='),\n    @('22+65=', '69-5='),\n    @('18-1=', '35-8='),\n    @('6+29=', '39-17='),\n    @('14-11=', '75-43='),\n    @('76-1=', '86-47='),\n    @('14+52=', '88-51='),\n    @('43-12=', '38+5='),\n    @('51+36=', '32+17='),\n    @('4+85=', '93-79='),\n    @('95-59=', '83-78='),\n    @('73-21=', '95-37='),\n    @('65+11=', '67+5='),\n    @('14+15=', '38-21='),\n    @('96-2=', '22+46='),\n    @('92-45=', '61-41='),\n    @('44+29=', '82-15='),\n    @('56+35=', '14+44='),\n    @('76+2=', '69-64='),\n    @('39-24=', '44+5='),\n    @('71-31=', '46+53='),\n    @('21+15=', '45-21='),\n    @('15+26=', '55+29='),\n    @('64-43=', '63+15='),\n    @('31+8=', '25+9='),\n    @('99-87=', '86-35='),\n    @('93+6=', '66-14='),\n    @('67+31=', '79+9='),\n    @('12+36=', '89-79='),\n    @('83-62=', '31+33='),\n    @('80-17=', '32+19='),\n    @('70-44=', '79-56='),\n    @('59-22=', '79-1='),\n    @('49-13=', '77+14='),\n    @('12+65=', '83-9='),\n    @('47+43=', '42+47='),\n    @('79-33=', '37+5='),\n    @('33-2=', '6+66='),\n    @('56-36=', '17-1='),\n    @('67-53=', '3+84='),\n    @('36-23=', '35+16='),\n    @('64+31=', '85-56='),\n    @('1+57=', '74-59='),\n    @('58-22=', '10+88='),\n    @('39-4=', '36-31='),\n    @('29-18=', '88-68='),\n    @('50-32=', '34-27='),\n    @('38-38=', '80-34='),\n    @('12+29=', '45-19='),\n    @('83-44=', '51+43='),\n    @('65-29=', '47-47='),\n    @('96-56=', '98-86='),\n    @('65-36=', '81-50='),\n    @('33-21=', '83-21='),\n    @('25+49=', '87-25='),\n    @('83-52=', '82-68='),\n    @('43+0=', '26+1='),\n    @('38-4=', '40-27='),\n    @('99-31=', '92-28='),\n    @('24-22=', '38-35='),\n    @('79-15=', '30-5='),\n    @('42+44=', '2+94='),\n    @('69+22=', '32+12='),\n    @('12+77=', '66-19='),\n    @('75-54=', '85-18='),\n    @('52-35=', '13+51='),\n    @('99-98=', '65+32='),\n    @('84-2=', '41+52='),\n    @('26+11=', '63+12='),\n    @('73+0=', '56-44='),\n    @('53-1=', '46+20='),\n    @('76-72=', '47-19='),\n    @('95-42=', '57+24='),\n    @('6+91=', '81-40='),\n    @('60+27=', '21+25='),\n    @('41-9=', '9+72='),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}"}
